$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.271.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.929.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7221'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -10.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3277'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06818'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8049'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07955'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.929.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.423'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '261.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.269.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007942'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.817'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.182.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9989'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.892'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.698'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1349'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.287'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.548'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.399'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.201'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05089'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.201'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7434'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.725'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01940'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.582'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("E42").Value = '  -5.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.010'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8350'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.727'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.288'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4118'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.483'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.93%  '
